$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 13:05"

# Row 13: Iran - refreshed case counts
$ws.Range("B13").Value = 114533
$ws.Range("C13").Value = 1808
$ws.Range("D13").Value = 90539
$ws.Range("E13").Value = 17140
$ws.Range("F13").Value = 2758
$ws.Range("G13").Value = 71
$ws.Range("H13").Value = 6854

# Rows 26-30: countries reshuffled due to re-sort by total cases, plus refreshed data
# Row 26 -> Catar
$ws.Range("A26").Value = "Catar"
$ws.Range("B26").Value = 28272
$ws.Range("C26").Value = 1733
$ws.Range("D26").Value = 3356
$ws.Range("E26").Value = 24902
$ws.Range("F26").Value = 72
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 14

# Row 27 -> Portugal
$ws.Range("A27").Value = "Portugal"
$ws.Range("B27").Value = 28132
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 3182
$ws.Range("E27").Value = 23775
$ws.Range("F27").Value = 103
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 1175

# Row 28 -> Suecia
$ws.Range("A28").Value = "Suecia"
$ws.Range("B28").Value = 27909
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 4971
$ws.Range("E28").Value = 19478
$ws.Range("F28").Value = 351
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 3460

# Row 29 -> Bielorrusia
$ws.Range("A29").Value = "Bielorrusia"
$ws.Range("B29").Value = 26772
$ws.Range("C29").Value = 947
$ws.Range("D29").Value = 8168
$ws.Range("E29").Value = 18453
$ws.Range("F29").Value = 92
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 151

# Row 30 -> Singapur
$ws.Range("A30").Value = "Singapur"
$ws.Range("B30").Value = 26098
$ws.Range("C30").Value = 752
$ws.Range("D30").Value = 4809
$ws.Range("E30").Value = 21268
$ws.Range("F30").Value = 19
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 21

# Row 31 (Irlanda) - unchanged, no edit needed

# Row 32: Emiratos Arabes Unidos - refreshed data only, country unchanged
$ws.Range("B32").Value = 21084
$ws.Range("C32").Value = 698
$ws.Range("D32").Value = 6930
$ws.Range("E32").Value = 13946
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 208

# Row 107: Libano - refreshed data
$ws.Range("B107").Value = 886
$ws.Range("C107").Value = 8
$ws.Range("E107").Value = 624

# Row 137: Montenegro - refreshed data
$ws.Range("D137").Value = 309
$ws.Range("E137").Value = 6
